$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Bugs"
$ws.Range("J2").Value = "20/06/2023 - 25/06/2023"

$ws.Range("J3").Select()
